$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.09"
$ws.Range("E2").Value = "'0.19%"
$ws.Range("D3").Value = "'40.74"
$ws.Range("E3").Value = "'-0.65%"
$ws.Range("D4").Value = "'5.126"
$ws.Range("E4").Value = "'1.51%"
$ws.Range("E5").Value = "'-0.07%"
$ws.Range("E6").Value = "'1.03%"
$ws.Range("D7").Value = "'1.610"
$ws.Range("E7").Value = "'0.73%"
$ws.Range("D8").Value = "'2.465"
$ws.Range("E8").Value = "'1.89%"
$ws.Range("D9").Value = "'0.9090"
$ws.Range("E9").Value = "'0.38%"
$ws.Range("E10").Value = "'27.99%"
$ws.Range("E11").Value = "'2.32%"
$ws.Range("D12").Value = "'0.09073"
$ws.Range("E12").Value = "'-1.09%"
$ws.Range("D13").Value = "'0.04320"
$ws.Range("E13").Value = "'-2.22%"
$ws.Range("E14").Value = "'-0.80%"
$ws.Range("E15").Value = "'-0.31%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005889"
$ws.Range("E16").Value = "'1.17%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.347"
$ws.Range("E17").Value = "'-0.54%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3314"
$ws.Range("E18").Value = "'0.48%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'6.988"
$ws.Range("E19").Value = "'3.33%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1394"
$ws.Range("E20").Value = "'3.27%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2706"
$ws.Range("E21").Value = "'-4.82%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04048"
$ws.Range("E22").Value = "'-2.77%"
$ws.Range("D23").Value = "'0.001275"
$ws.Range("E23").Value = "'4.95%"
$ws.Range("E24").Value = "'-0.78%"
$ws.Range("D25").Value = "'0.0001275"
$ws.Range("E25").Value = "'-1.88%"
$ws.Range("E26").Value = "'24.67%"
$ws.Range("D38").Value = "'0.02421"
$ws.Range("E38").Value = "'0.21%"
$ws.Range("D39").Value = "'0.05227"
$ws.Range("E39").Value = "'1.57%"
$ws.Range("D40").Value = "'0.007836"
$ws.Range("E40").Value = "'0.07%"
$ws.Range("D41").Value = "'0.1300"
$ws.Range("E41").Value = "'-0.39%"
$ws.Range("D42").Value = "'0.006813"
$ws.Range("E42").Value = "'-3.64%"
$ws.Range("D43").Value = "'0.001903"
$ws.Range("E43").Value = "'-2.30%"
$ws.Range("D44").Value = "'0.007409"
$ws.Range("E44").Value = "'-10.28%"
$ws.Range("D45").Value = "'0.3337"
$ws.Range("E45").Value = "'9.70%"
$ws.Range("D46").Value = "'0.00006911"
$ws.Range("E46").Value = "'8.24%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.31%"
$ws.Range("D48").Value = "'0.1092"
$ws.Range("E48").Value = "'1,696.75%"
$ws.Range("D50").Value = "'0.00002108"
$ws.Range("E50").Value = "'0.31%"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("E51").Value = "'0.31%"
